# quickfix - strategy headers not being written correctly
#
# The three "ORDER STATISTIC - ..." section headers (A2, I2, Q2) were all
# being written with the same text ("ORDER STATISTIC - QUICK SORT"). Each
# column block actually reports a different selection strategy, so fix the
# header text for the iterative-selection block (I2) and the quick-select
# block (Q2). A2 (quick sort) is already correct and is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = "ORDER STATISTIC - ITERATIVE SELECTION"
$ws.Range("Q2").Value = "ORDER STATISTIC - QUICK SELECT"

# Move the active selection to Q3, matching the saved cursor position.
$ws.Range("Q3").Select()
